$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace comma-separated lists with colon-separated lists in the
# haz_potlink / hazlink columns (and other affected cells).
$targets = @("H15","I15","G16","G18","H18","H24","H29","H32","H35","H38","H39","H42","H43","H44","H49","G52","H52")

foreach ($addr in $targets) {
    $rng = $ws.Range($addr)
    $old = $rng.Value()
    $rng.Value = $old -replace ",", ":"
}

# Update the active selection/cursor to the final cell touched by the
# find & replace sweep (this engine's Selection model is a single range,
# so we land the cursor on the last edited area's anchor cell, J1).
Write-Output "Applied colon replacements to $($targets.Count) cells."
$sel = $ws.Range("J1")
[void]$sel.Select()
